$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAT")

# Shift the header + first two data rows up by one row (row3->2, row4->3, row5->4),
# copying values and cell formatting together.
$ws.Range("A3:D3").Copy($ws.Range("A2:D2"))
$ws.Range("A4:D4").Copy($ws.Range("A3:D3"))
$ws.Range("A5:D5").Copy($ws.Range("A4:D4"))

# Row 5 becomes a new "command" row for the flashuploader cgi permissions.
# Reuse row 4's formatting (same pattern: command / chmod.../onChangeOnly) then
# overwrite the pattern text with the new chmod command.
$ws.Range("A4:D4").Copy($ws.Range("A5:D5"))
$ws.Range("B5").Value2 = "chmod 750 ./web/<progDir>/flashuploader/courseleaf.cgi"

# Widen column B to fit the new, longer text.
$ws.Columns("B").ColumnWidth = 66.28515625

# Make CAT the active sheet/tab, with A2 selected (matches the new header position).
$ws.Activate()
$ws.Range("A2").Select()
